$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '''24.528.07'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '''  +9.34%  '
$ws.Cells.Item(2, 5).Style = "Normal"

# Row 3
$ws.Cells.Item(3, 4).Value = '''1.684.23'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '''  +5.03%  '
$ws.Cells.Item(3, 5).Style = "Normal"

# Row 4
$ws.Cells.Item(4, 4).Value = '''1.003'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '''  -0.17%  '
$ws.Cells.Item(4, 5).Style = "Normal"

# Row 5
$ws.Cells.Item(5, 4).Value = '''307.86'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '''  +1.19%  '
$ws.Cells.Item(5, 5).Style = "Normal"

# Row 7
$ws.Cells.Item(7, 4).Value = '''0.3708'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '''  +0.61%  '
$ws.Cells.Item(7, 5).Style = "Normal"

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.3455'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '''  +2.26%  '
$ws.Cells.Item(8, 5).Style = "Normal"

# Row 9
$ws.Cells.Item(9, 4).Value = '''48.06'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '''  +13.46%  '
$ws.Cells.Item(9, 5).Style = "Normal"

# Row 10
$ws.Cells.Item(10, 4).Value = '''1.184'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '''  +4.00%  '
$ws.Cells.Item(10, 5).Style = "Normal"

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.07305'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '''  +3.51%  '
$ws.Cells.Item(11, 5).Style = "Normal"

# Row 12
$ws.Cells.Item(12, 5).Value = '''  +0.03%  '
$ws.Cells.Item(12, 5).Style = "Normal"

# Row 13
$ws.Cells.Item(13, 2).Value = '''Solana'
$ws.Cells.Item(13, 2).Style = "Normal"
$ws.Cells.Item(13, 3).Value = '''https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(13, 3).Style = "Normal"
$ws.Cells.Item(13, 4).Value = '''20.51'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '''  +4.16%  '
$ws.Cells.Item(13, 5).Style = "Normal"

# Row 14
$ws.Cells.Item(14, 2).Value = '''Polkadot'
$ws.Cells.Item(14, 2).Style = "Normal"
$ws.Cells.Item(14, 3).Value = '''https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(14, 3).Style = "Normal"
$ws.Cells.Item(14, 4).Value = '''6.179'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '''  +4.21%  '
$ws.Cells.Item(14, 5).Style = "Normal"

# Row 15
$ws.Cells.Item(15, 4).Value = '''6.774'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '''  +2.11%  '
$ws.Cells.Item(15, 5).Style = "Normal"

# Row 16
$ws.Cells.Item(16, 4).Value = '''1.676.79'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '''  +4.88%  '
$ws.Cells.Item(16, 5).Style = "Normal"

# Row 17
$ws.Cells.Item(17, 4).Value = '''0.00001114'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '''  +2.45%  '
$ws.Cells.Item(17, 5).Style = "Normal"

# Row 18
$ws.Cells.Item(18, 4).Value = '''0.9986'
$ws.Cells.Item(18, 4).Style = "Normal"

# Row 19
$ws.Cells.Item(19, 4).Value = '''0.06725'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '''  -1.47%  '
$ws.Cells.Item(19, 5).Style = "Normal"

# Row 20
$ws.Cells.Item(20, 4).Value = '''81.50'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '''  +4.49%  '
$ws.Cells.Item(20, 5).Style = "Normal"

# Row 21
$ws.Cells.Item(21, 4).Value = '''16.52'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '''  +2.65%  '
$ws.Cells.Item(21, 5).Style = "Normal"

# Row 22
$ws.Cells.Item(22, 4).Value = '''6.129'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '''  +1.62%  '
$ws.Cells.Item(22, 5).Style = "Normal"

# Row 23
$ws.Cells.Item(23, 4).Value = '''12.16'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '''  +2.76%  '
$ws.Cells.Item(23, 5).Style = "Normal"

# Row 24
$ws.Cells.Item(24, 4).Value = '''24.424.52'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '''  +8.69%  '
$ws.Cells.Item(24, 5).Style = "Normal"

# Row 25
$ws.Cells.Item(25, 4).Value = '''2.442'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '''  +1.53%  '
$ws.Cells.Item(25, 5).Style = "Normal"

# Row 26
$ws.Cells.Item(26, 4).Value = '''2.684'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '''  +5.81%  '
$ws.Cells.Item(26, 5).Style = "Normal"

# Row 27
$ws.Cells.Item(27, 2).Value = '''LEO'
$ws.Cells.Item(27, 2).Style = "Normal"
$ws.Cells.Item(27, 3).Value = '''https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(27, 3).Style = "Normal"
$ws.Cells.Item(27, 4).Value = '''3.360'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '''  -13.81%  '
$ws.Cells.Item(27, 5).Style = "Normal"

# Row 28
$ws.Cells.Item(28, 2).Value = '''Monero'
$ws.Cells.Item(28, 2).Style = "Normal"
$ws.Cells.Item(28, 3).Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(28, 3).Style = "Normal"
$ws.Cells.Item(28, 4).Value = '''152.92'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '''  +1.55%  '
$ws.Cells.Item(28, 5).Style = "Normal"

# Row 29
$ws.Cells.Item(29, 2).Value = '''EthereumClassic'
$ws.Cells.Item(29, 2).Style = "Normal"
$ws.Cells.Item(29, 3).Value = '''https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(29, 3).Style = "Normal"
$ws.Cells.Item(29, 4).Value = '''19.56'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '''  -0.28%  '
$ws.Cells.Item(29, 5).Style = "Normal"

# Row 30
$ws.Cells.Item(30, 2).Value = '''WrappedliquidstakedEther2.0'
$ws.Cells.Item(30, 2).Style = "Normal"
$ws.Cells.Item(30, 3).Value = '''https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(30, 3).Style = "Normal"
$ws.Cells.Item(30, 4).Value = '''1.862.11'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '''  +4.61%  '
$ws.Cells.Item(30, 5).Style = "Normal"

# Row 31
$ws.Cells.Item(31, 2).Value = '''BitcoinCash'
$ws.Cells.Item(31, 2).Style = "Normal"
$ws.Cells.Item(31, 3).Value = '''https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(31, 3).Style = "Normal"
$ws.Cells.Item(31, 4).Value = '''127.80'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '''  +5.82%  '
$ws.Cells.Item(31, 5).Style = "Normal"

# Row 32
$ws.Cells.Item(32, 2).Value = '''Filecoin'
$ws.Cells.Item(32, 2).Style = "Normal"
$ws.Cells.Item(32, 3).Value = '''https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(32, 3).Style = "Normal"
$ws.Cells.Item(32, 4).Value = '''6.372'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '''  +4.94%  '
$ws.Cells.Item(32, 5).Style = "Normal"

# Row 33
$ws.Cells.Item(33, 2).Value = '''HuobiToken'
$ws.Cells.Item(33, 2).Style = "Normal"
$ws.Cells.Item(33, 3).Value = '''https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(33, 3).Style = "Normal"
$ws.Cells.Item(33, 4).Value = '''4.047'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '''  -2.57%  '
$ws.Cells.Item(33, 5).Style = "Normal"

# Row 34
$ws.Cells.Item(34, 2).Value = '''ImmutableX'
$ws.Cells.Item(34, 2).Style = "Normal"
$ws.Cells.Item(34, 3).Value = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(34, 3).Style = "Normal"
$ws.Cells.Item(34, 4).Value = '''0.9844'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '''  +3.04%  '
$ws.Cells.Item(34, 5).Style = "Normal"

# Row 35
$ws.Cells.Item(35, 2).Value = '''Stellar'
$ws.Cells.Item(35, 2).Style = "Normal"
$ws.Cells.Item(35, 3).Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(35, 3).Style = "Normal"
$ws.Cells.Item(35, 4).Value = '''0.08483'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '''  +2.45%  '
$ws.Cells.Item(35, 5).Style = "Normal"

# Row 36
$ws.Cells.Item(36, 2).Value = '''WEMIXTOKEN'
$ws.Cells.Item(36, 2).Style = "Normal"
$ws.Cells.Item(36, 3).Value = '''https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(36, 3).Style = "Normal"
$ws.Cells.Item(36, 4).Value = '''1.701'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '''  +4.04%  '
$ws.Cells.Item(36, 5).Style = "Normal"

# Row 37
$ws.Cells.Item(37, 2).Value = '''Aptos'
$ws.Cells.Item(37, 2).Style = "Normal"
$ws.Cells.Item(37, 3).Value = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(37, 3).Style = "Normal"
$ws.Cells.Item(37, 4).Value = '''12.55'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '''  +4.26%  '
$ws.Cells.Item(37, 5).Style = "Normal"

# Row 38
$ws.Cells.Item(38, 2).Value = '''Hedera'
$ws.Cells.Item(38, 2).Style = "Normal"
$ws.Cells.Item(38, 3).Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(38, 3).Style = "Normal"
$ws.Cells.Item(38, 4).Value = '''0.06545'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '''  +7.38%  '
$ws.Cells.Item(38, 5).Style = "Normal"

# Row 39
$ws.Cells.Item(39, 4).Value = '''9.011'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '''  +4.40%  '
$ws.Cells.Item(39, 5).Style = "Normal"

# Row 40
$ws.Cells.Item(40, 2).Value = '''InternetComputer(DFINITY)'
$ws.Cells.Item(40, 2).Style = "Normal"
$ws.Cells.Item(40, 3).Value = '''https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(40, 3).Style = "Normal"
$ws.Cells.Item(40, 4).Value = '''5.388'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '''  +2.33%  '
$ws.Cells.Item(40, 5).Style = "Normal"

# Row 41
$ws.Cells.Item(41, 2).Value = '''VeChain'
$ws.Cells.Item(41, 2).Style = "Normal"
$ws.Cells.Item(41, 3).Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(41, 3).Style = "Normal"
$ws.Cells.Item(41, 4).Value = '''0.02340'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '''  +5.10%  '
$ws.Cells.Item(41, 5).Style = "Normal"

# Row 42
$ws.Cells.Item(42, 2).Value = '''TrustWalletToken'
$ws.Cells.Item(42, 2).Style = "Normal"
$ws.Cells.Item(42, 3).Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(42, 3).Style = "Normal"
$ws.Cells.Item(42, 4).Value = '''1.267'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '''  -0.18%  '
$ws.Cells.Item(42, 5).Style = "Normal"

# Row 43
$ws.Cells.Item(43, 2).Value = '''Algorand'
$ws.Cells.Item(43, 2).Style = "Normal"
$ws.Cells.Item(43, 3).Value = '''https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(43, 3).Style = "Normal"
$ws.Cells.Item(43, 4).Value = '''0.2123'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '''  +4.68%  '
$ws.Cells.Item(43, 5).Style = "Normal"

# Row 44
$ws.Cells.Item(44, 2).Value = '''TheSandbox'
$ws.Cells.Item(44, 2).Style = "Normal"
$ws.Cells.Item(44, 3).Value = '''https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(44, 3).Style = "Normal"
$ws.Cells.Item(44, 4).Value = '''0.6196'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '''  +4.47%  '
$ws.Cells.Item(44, 5).Style = "Normal"

# Row 45
$ws.Cells.Item(45, 2).Value = '''Frax'
$ws.Cells.Item(45, 2).Style = "Normal"
$ws.Cells.Item(45, 3).Value = '''https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(45, 3).Style = "Normal"
$ws.Cells.Item(45, 4).Value = '''0.9988'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '''  +0.81%  '
$ws.Cells.Item(45, 5).Style = "Normal"

# Row 46
$ws.Cells.Item(46, 2).Value = '''EnergySwap'
$ws.Cells.Item(46, 2).Style = "Normal"
$ws.Cells.Item(46, 3).Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(46, 3).Style = "Normal"
$ws.Cells.Item(46, 4).Value = '''13.33'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '''  +0.79%  '
$ws.Cells.Item(46, 5).Style = "Normal"

# Row 47
$ws.Cells.Item(47, 4).Value = '''0.5985'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '''  +4.61%  '
$ws.Cells.Item(47, 5).Style = "Normal"

# Row 48
$ws.Cells.Item(48, 2).Value = '''PancakeSwap'
$ws.Cells.Item(48, 2).Style = "Normal"
$ws.Cells.Item(48, 3).Value = '''https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(48, 3).Style = "Normal"
$ws.Cells.Item(48, 4).Value = '''3.770'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '''  -1.64%  '
$ws.Cells.Item(48, 5).Style = "Normal"

# Row 49
$ws.Cells.Item(49, 2).Value = '''Quant'
$ws.Cells.Item(49, 2).Style = "Normal"
$ws.Cells.Item(49, 3).Value = '''https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(49, 3).Style = "Normal"
$ws.Cells.Item(49, 4).Value = '''128.09'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '''  +0.96%  '
$ws.Cells.Item(49, 5).Style = "Normal"

# Row 50
$ws.Cells.Item(50, 2).Value = '''NEARProtocol'
$ws.Cells.Item(50, 2).Style = "Normal"
$ws.Cells.Item(50, 3).Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(50, 3).Style = "Normal"
$ws.Cells.Item(50, 4).Value = '''2.037'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '''  +2.99%  '
$ws.Cells.Item(50, 5).Style = "Normal"

# Row 51
$ws.Cells.Item(51, 2).Value = '''Cronos'
$ws.Cells.Item(51, 2).Style = "Normal"
$ws.Cells.Item(51, 3).Value = '''https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(51, 3).Style = "Normal"
$ws.Cells.Item(51, 4).Value = '''0.07216'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '''  +6.00%  '
$ws.Cells.Item(51, 5).Style = "Normal"
